{"js": "// This document is a title line (\"YYYY-MM-DD Weekday\") followed by a\n// 5-column table of 100 two-term arithmetic problems (\"a+b=c\" / \"a-b=c\"),\n// one expression per cell. The edit replaces the title's date and every\n// table-cell expression with a new value, each taken from the same\n// position (one <w:t> run per entry, read in document order: the title\n// first, then row-by-row/left-to-right through the table).\n//\n// `pairs` below holds one [oldText, newText] tuple per run, in that same\n// document order. A couple of positions happen to share the exact same\n// oldText (\"6+87=93\" appears twice, with two different replacements), so\n// a plain find-and-replace-all over the whole document would not be able\n// to tell the two occurrences apart. Instead, for each distinct oldText\n// we gather every place it occurs (in document order) and pair those\n// occurrences up, in order, with the list of intended replacements for\n// that text.\nconst pairs = [\n  [\"2025-12-29 Monday\", \"2025-12-30 Tuesday\"],\n  [\"82-59=23\", \"16+6=22\"],\n  [\"61-8=53\", \"82-46=36\"],\n  [\"29+46=75\", \"74-27=47\"],\n  [\"44-38=6\", \"66-9=57\"],\n  [\"28+16=44\", \"50-28=22\"],\n  [\"67+18=85\", \"48+13=61\"],\n  [\"79+19=98\", \"61-26=35\"],\n  [\"40-37=3\", \"70-51=19\"],\n  [\"6+87=93\", \"32-26=6\"],\n  [\"30-15=15\", \"74-68=6\"],\n  [\"78+3=81\", \"61-17=44\"],\n  [\"63-34=29\", \"29+24=53\"],\n  [\"6+36=42\", \"28+17=45\"],\n  [\"53-6=47\", \"82-65=17\"],\n  [\"45+39=84\", \"55-16=39\"],\n  [\"45-28=17\", \"24+38=62\"],\n  [\"66+29=95\", \"31-15=16\"],\n  [\"23-8=15\", \"18+43=61\"],\n  [\"19+59=78\", \"24-6=18\"],\n  [\"60-42=18\", \"58+16=74\"],\n  [\"24+8=32\", \"22-19=3\"],\n  [\"49+28=77\", \"16+16=32\"],\n  [\"61-45=16\", \"27+17=44\"],\n  [\"60-47=13\", \"91-89=2\"],\n  [\"28+46=74\", \"68+15=83\"],\n  [\"34+48=82\", \"71-59=12\"],\n  [\"80-46=34\", \"93-45=48\"],\n  [\"71-15=56\", \"96-48=48\"],\n  [\"93-47=46\", \"14+77=91\"],\n  [\"39+17=56\", \"7+35=42\"],\n  [\"47+44=91\", \"24+67=91\"],\n  [\"87-38=49\", \"91-53=38\"],\n  [\"9+85=94\", \"65+6=71\"],\n  [\"30-22=8\", \"12-3=9\"],\n  [\"41-3=38\", \"80-25=55\"],\n  [\"71-57=14\", \"84+9=93\"],\n  [\"81-22=59\", \"90-12=78\"],\n  [\"26-18=8\", \"94-36=58\"],\n  [\"8+85=93\", \"49+33=82\"],\n  [\"8+25=33\", \"65+9=74\"],\n  [\"53+9=62\", \"41-18=23\"],\n  [\"93-86=7\", \"74-26=48\"],\n  [\"92-59=33\", \"85+6=91\"],\n  [\"90-22=68\", \"39+49=88\"],\n  [\"9+48=57\", \"22-9=13\"],\n  [\"59+4=63\", \"72-4=68\"],\n  [\"53-45=8\", \"38+53=91\"],\n  [\"72-64=8\", \"97-39=58\"],\n  [\"80-42=38\", \"62-19=43\"],\n  [\"16+15=31\", \"33-19=14\"],\n  [\"84-37=47\", \"91-36=55\"],\n  [\"51-48=3\", \"81-42=39\"],\n  [\"56-48=8\", \"58+17=75\"],\n  [\"93-17=76\", \"35+16=51\"],\n  [\"96-7=89\", \"78+9=87\"],\n  [\"72-14=58\", \"8+36=44\"],\n  [\"27+45=72\", \"87-8=79\"],\n  [\"90-66=24\", \"76+8=84\"],\n  [\"15+7=22\", \"9+57=66\"],\n  [\"6+87=93\", \"61-36=25\"],\n  [\"75-67=8\", \"51-2=49\"],\n  [\"36+5=41\", \"46+7=53\"],\n  [\"90-64=26\", \"72-33=39\"],\n  [\"37+58=95\", \"90-51=39\"],\n  [\"82-36=46\", \"85-28=57\"],\n  [\"72-13=59\", \"94-17=77\"],\n  [\"18+3=21\", \"85+9=94\"],\n  [\"18+45=63\", \"31-28=3\"],\n  [\"55-26=29\", \"40-13=27\"],\n  [\"63+28=91\", \"77-49=28\"],\n  [\"17+9=26\", \"92-16=76\"],\n  [\"6+55=61\", \"67+28=95\"],\n  [\"7+4=11\", \"40-38=2\"],\n  [\"8+68=76\", \"5+38=43\"],\n  [\"16+9=25\", \"67+25=92\"],\n  [\"26+65=91\", \"80-33=47\"],\n  [\"92-77=15\", \"92-18=74\"],\n  [\"75-17=58\", \"19+68=87\"],\n  [\"7+9=16\", \"75-9=66\"],\n  [\"54+28=82\", \"25+7=32\"],\n  [\"84-55=29\", \"90-28=62\"],\n  [\"81-66=15\", \"19+46=65\"],\n  [\"53+18=71\", \"9+84=93\"],\n  [\"48+5=53\", \"28+35=63\"],\n  [\"94-38=56\", \"16+65=81\"],\n  [\"12+59=71\", \"52-24=28\"],\n  [\"85-47=38\", \"71-28=43\"],\n  [\"84+8=92\", \"67+17=84\"],\n  [\"62-4=58\", \"77+17=94\"],\n  [\"53-36=17\", \"97-38=59\"],\n  [\"60-43=17\", \"25-17=8\"],\n  [\"74-29=45\", \"8+17=25\"],\n  [\"57+14=71\", \"65-49=16\"],\n  [\"14-6=8\", \"9+65=74\"],\n  [\"6+28=34\", \"73-46=27\"],\n  [\"19+73=92\", \"85-37=48\"],\n  [\"51-32=19\", \"90-56=34\"],\n  [\"65+29=94\", \"45+48=93\"],\n  [\"57+28=85\", \"91-48=43\"],\n  [\"81-13=68\", \"18+36=54\"]\n];\n\nconst byOld = new Map();\nfor (const [oldText, newText] of pairs) {\n  if (!byOld.has(oldText)) byOld.set(oldText, []);\n  byOld.get(oldText).push(newText);\n}\n\nfor (const [oldText, newTexts] of byOld) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== newTexts.length) {\n    throw new Error(\n      `Expected ${newTexts.length} match(es) for \"${oldText}\" but found ${results.items.length}`\n    );\n  }\n\n  // results.items are returned in document order, matching the order of\n  // newTexts collected above.\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newTexts[i], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Ordered list of (old, new) text pairs, one per run in document order\n# (title paragraph's date, then each of the 100 table-cell arithmetic\n# expressions). One old value (\"6+87=93\") occurs twice in the original\n# document with two different replacements, so each pair is resolved by\n# advancing a cursor through the document and replacing only the first\n# remaining match for each pair's old text -- never a blanket\n# find-replace-all, which would not be able to distinguish the two\n# differing occurrences of \"6+87=93\".\n$pairs = @(\n  @(\"2025-12-29 Monday\", \"2025-12-30 Tuesday\"),\n  @(\"82-59=23\", \"16+6=22\"),\n  @(\"61-8=53\", \"82-46=36\"),\n  @(\"29+46=75\", \"74-27=47\"),\n  @(\"44-38=6\", \"66-9=57\"),\n  @(\"28+16=44\", \"50-28=22\"),\n  @(\"67+18=85\", \"48+13=61\"),\n  @(\"79+19=98\", \"61-26=35\"),\n  @(\"40-37=3\", \"70-51=19\"),\n  @(\"6+87=93\", \"32-26=6\"),\n  @(\"30-15=15\", \"74-68=6\"),\n  @(\"78+3=81\", \"61-17=44\"),\n  @(\"63-34=29\", \"29+24=53\"),\n  @(\"6+36=42\", \"28+17=45\"),\n  @(\"53-6=47\", \"82-65=17\"),\n  @(\"45+39=84\", \"55-16=39\"),\n  @(\"45-28=17\", \"24+38=62\"),\n  @(\"66+29=95\", \"31-15=16\"),\n  @(\"23-8=15\", \"18+43=61\"),\n  @(\"19+59=78\", \"24-6=18\"),\n  @(\"60-42=18\", \"58+16=74\"),\n  @(\"24+8=32\", \"22-19=3\"),\n  @(\"49+28=77\", \"16+16=32\"),\n  @(\"61-45=16\", \"27+17=44\"),\n  @(\"60-47=13\", \"91-89=2\"),\n  @(\"28+46=74\", \"68+15=83\"),\n  @(\"34+48=82\", \"71-59=12\"),\n  @(\"80-46=34\", \"93-45=48\"),\n  @(\"71-15=56\", \"96-48=48\"),\n  @(\"93-47=46\", \"14+77=91\"),\n  @(\"39+17=56\", \"7+35=42\"),\n  @(\"47+44=91\", \"24+67=91\"),\n  @(\"87-38=49\", \"91-53=38\"),\n  @(\"9+85=94\", \"65+6=71\"),\n  @(\"30-22=8\", \"12-3=9\"),\n  @(\"41-3=38\", \"80-25=55\"),\n  @(\"71-57=14\", \"84+9=93\"),\n  @(\"81-22=59\", \"90-12=78\"),\n  @(\"26-18=8\", \"94-36=58\"),\n  @(\"8+85=93\", \"49+33=82\"),\n  @(\"8+25=33\", \"65+9=74\"),\n  @(\"53+9=62\", \"41-18=23\"),\n  @(\"93-86=7\", \"74-26=48\"),\n  @(\"92-59=33\", \"85+6=91\"),\n  @(\"90-22=68\", \"39+49=88\"),\n  @(\"9+48=57\", \"22-9=13\"),\n  @(\"59+4=63\", \"72-4=68\"),\n  @(\"53-45=8\", \"38+53=91\"),\n  @(\"72-64=8\", \"97-39=58\"),\n  @(\"80-42=38\", \"62-19=43\"),\n  @(\"16+15=31\", \"33-19=14\"),\n  @(\"84-37=47\", \"91-36=55\"),\n  @(\"51-48=3\", \"81-42=39\"),\n  @(\"56-48=8\", \"58+17=75\"),\n  @(\"93-17=76\", \"35+16=51\"),\n  @(\"96-7=89\", \"78+9=87\"),\n  @(\"72-14=58\", \"8+36=44\"),\n  @(\"27+45=72\", \"87-8=79\"),\n  @(\"90-66=24\", \"76+8=84\"),\n  @(\"15+7=22\", \"9+57=66\"),\n  @(\"6+87=93\", \"61-36=25\"),\n  @(\"75-67=8\", \"51-2=49\"),\n  @(\"36+5=41\", \"46+7=53\"),\n  @(\"90-64=26\", \"72-33=39\"),\n  @(\"37+58=95\", \"90-51=39\"),\n  @(\"82-36=46\", \"85-28=57\"),\n  @(\"72-13=59\", \"94-17=77\"),\n  @(\"18+3=21\", \"85+9=94\"),\n  @(\"18+45=63\", \"31-28=3\"),\n  @(\"55-26=29\", \"40-13=27\"),\n  @(\"63+28=91\", \"77-49=28\"),\n  @(\"17+9=26\", \"92-16=76\"),\n  @(\"6+55=61\", \"67+28=95\"),\n  @(\"7+4=11\", \"40-38=2\"),\n  @(\"8+68=76\", \"5+38=43\"),\n  @(\"16+9=25\", \"67+25=92\"),\n  @(\"26+65=91\", \"80-33=47\"),\n  @(\"92-77=15\", \"92-18=74\"),\n  @(\"75-17=58\", \"19+68=87\"),\n  @(\"7+9=16\", \"75-9=66\"),\n  @(\"54+28=82\", \"25+7=32\"),\n  @(\"84-55=29\", \"90-28=62\"),\n  @(\"81-66=15\", \"19+46=65\"),\n  @(\"53+18=71\", \"9+84=93\"),\n  @(\"48+5=53\", \"28+35=63\"),\n  @(\"94-38=56\", \"16+65=81\"),\n  @(\"12+59=71\", \"52-24=28\"),\n  @(\"85-47=38\", \"71-28=43\"),\n  @(\"84+8=92\", \"67+17=84\"),\n  @(\"62-4=58\", \"77+17=94\"),\n  @(\"53-36=17\", \"97-38=59\"),\n  @(\"60-43=17\", \"25-17=8\"),\n  @(\"74-29=45\", \"8+17=25\"),\n  @(\"57+14=71\", \"65-49=16\"),\n  @(\"14-6=8\", \"9+65=74\"),\n  @(\"6+28=34\", \"73-46=27\"),\n  @(\"19+73=92\", \"85-37=48\"),\n  @(\"51-32=19\", \"90-56=34\"),\n  @(\"65+29=94\", \"45+48=93\"),\n  @(\"57+28=85\", \"91-48=43\"),\n  @(\"81-13=68\", \"18+36=54\")\n)\n\n$d = $word.ActiveDocument\n$cursorStart = 0\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  # Search only the remaining (not-yet-processed) tail of the document so\n  # that repeated old values resolve to the correct successive new value.\n  $searchRange = $d.Range($cursorStart, $d.Content.End)\n  $searchRange.Find.ClearFormatting()\n  $found = $searchRange.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n\n  if (-not $found) {\n    throw \"Could not find expected text '$oldText' while applying replacements\"\n  }\n\n  # Resume the next search right after this replacement.\n  $cursorStart = $searchRange.End\n}\n\n"}
